$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename first sheet: "RiskList" -> "RiskList v1"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RiskList v1"

# ---------------------------------------------------------------------------
# 2. RiskList v2 sheet: add a "Stav" (Status) column with "Aktivní" (Active)
#    for every data row, then sort the table by Priority desc, Probability desc.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("I1").Value = "Stav"
for ($r = 2; $r -le 11; $r++) {
    $ws2.Range("I$r").Value = "Aktivní"
}

# New column inherits the formatting of column H (same row banding/border).
$ws2.Range("H1:H11").Copy()
$ws2.Range("I1:I11").PasteSpecial(-4122)   # xlPasteFormats

# Row for ID=6 (currently row 7): make the rating cells match the rest of the row.
$ws2.Range("A7").Copy()
$ws2.Range("E7:F7").PasteSpecial(-4122)
$ws2.Range("H7").PasteSpecial(-4122)

# Row for ID=8 (currently row 9): give the rating cells (and new Stav cell) a
# bold box border while keeping the row's existing highlight colour.
$ws2.Range("A2").Copy()
$ws2.Range("E9:F9").PasteSpecial(-4122)
$ws2.Range("H9:I9").PasteSpecial(-4122)
$fillColor = $ws2.Range("G9").Interior.Color
$ws2.Range("E9:F9").Interior.Color = $fillColor
$ws2.Range("H9:I9").Interior.Color = $fillColor

$excel.CutCopyMode = 0

# Sort the table by Priority (column E) descending, then Probability
# (column F) descending - matching the dialog settings used in the workbook.
$sortRange = $ws2.Range("A1:I11")
$key1 = $ws2.Range("E2:E11")
$key2 = $ws2.Range("F2:F11")
$sortRange.Sort($key1, 2, $key2, $null, 2, $null, $null, 1)
